# -----------------------------------------------------------------------------
# Commit: "update scripts wuth new tpm"
#
# The NATMI ligand-receptor edge table (Ntrk3 -> Ptprs) was regenerated against
# an updated TPM expression matrix. All computed metric columns E:T (ligand/
# receptor expression, detection rate, specificity scores and edge weights) are
# refreshed with the newly-computed values for every Sending/Target cluster pair
# (rows 2-13). Columns A:D (cluster/gene identifiers) are untouched.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order for the range E:T being written on every data row:
# E  Ligand-expressing cells
# F  Ligand detection rate
# G  Ligand average expression value
# H  Ligand total expression value
# I  Ligand derived specificity of average expression value
# J  Ligand derived specificity of total expression value
# K  Receptor-expressing cells
# L  Receptor detection rate
# M  Receptor average expression value
# N  Receptor total expression value
# O  Receptor derived specificity of average expression value
# P  Receptor derived specificity of total expression value
# Q  Edge average expression weight
# R  Edge total expression weight
# S  Edge average expression derived specificity
# T  Edge total expression derived specificity

# Row 2: ECs -> ECs
$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 2
$row2[0,1] = 0.6666666666666666
$row2[0,2] = 0.1487266666666667
$row2[0,3] = 0.44618
$row2[0,4] = 0.1109321277273384
$row2[0,5] = 0.1109321277273384
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 3.556762333333333
$row2[0,9] = 10.670287
$row2[0,10] = 0.04280930450251701
$row2[0,11] = 0.04280930450251701
$row2[0,12] = 0.5289854059622222
$row2[0,13] = 4.76086865366
$row2[0,14] = 0.004748927234991739
$row2[0,15] = 0.004748927234991739
$ws.Range("E2:T2").Value = $row2

# Row 3: ECs -> FAPs
$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 2
$row3[0,1] = 0.6666666666666666
$row3[0,2] = 0.1487266666666667
$row3[0,3] = 0.44618
$row3[0,4] = 0.1109321277273384
$row3[0,5] = 0.1109321277273384
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 47.24901333333333
$row3[0,9] = 141.74704
$row3[0,10] = 0.5686906263805706
$row3[0,11] = 0.5686906263805704
$row3[0,12] = 7.027188256355555
$row3[0,13] = 63.2446943072
$row3[0,14] = 0.06308606120298951
$row3[0,15] = 0.06308606120298951
$ws.Range("E3:T3").Value = $row3

# Row 4: ECs -> MuSCs
$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 2
$row4[0,1] = 0.6666666666666666
$row4[0,2] = 0.1487266666666667
$row4[0,3] = 0.44618
$row4[0,4] = 0.1109321277273384
$row4[0,5] = 0.1109321277273384
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 24.53173066666666
$row4[0,9] = 73.595192
$row4[0,10] = 0.2952646900921413
$row4[0,11] = 0.2952646900921412
$row4[0,12] = 3.648522529617777
$row4[0,13] = 32.83670276656
$row4[0,14] = 0.0327543403146744
$row4[0,15] = 0.0327543403146744
$ws.Range("E4:T4").Value = $row4

# Row 5: ECs -> Resolving-Mac
$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 2
$row5[0,1] = 0.6666666666666666
$row5[0,2] = 0.1487266666666667
$row5[0,3] = 0.44618
$row5[0,4] = 0.1109321277273384
$row5[0,5] = 0.1109321277273384
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 7.746355333333334
$row5[0,9] = 23.239066
$row5[0,10] = 0.09323537902477132
$row5[0,11] = 0.0932353790247713
$row5[0,12] = 1.152089607542222
$row5[0,13] = 10.36880646788
$row5[0,14] = 0.01034279897468274
$row5[0,15] = 0.01034279897468274
$ws.Range("E5:T5").Value = $row5

# Row 6: FAPs -> ECs
$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 0.07403633333333333
$row6[0,3] = 0.222109
$row6[0,4] = 0.05522216136400421
$row6[0,5] = 0.05522216136400422
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 3.556762333333333
$row6[0,9] = 10.670287
$row6[0,10] = 0.04280930450251701
$row6[0,11] = 0.04280930450251701
$row6[0,12] = 0.2633296416981111
$row6[0,13] = 2.369966775283
$row6[0,14] = 0.002364022321118786
$row6[0,15] = 0.002364022321118786
$ws.Range("E6:T6").Value = $row6

# Row 7: FAPs -> FAPs
$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 0.07403633333333333
$row7[0,3] = 0.222109
$row7[0,4] = 0.05522216136400421
$row7[0,5] = 0.05522216136400422
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 47.24901333333333
$row7[0,9] = 141.74704
$row7[0,10] = 0.5686906263805706
$row7[0,11] = 0.5686906263805704
$row7[0,12] = 3.498143700817777
$row7[0,13] = 31.48329330736
$row7[0,14] = 0.0314043255361845
$row7[0,15] = 0.0314043255361845
$ws.Range("E7:T7").Value = $row7

# Row 8: FAPs -> MuSCs
$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 0.07403633333333333
$row8[0,3] = 0.222109
$row8[0,4] = 0.05522216136400421
$row8[0,5] = 0.05522216136400422
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 24.53173066666666
$row8[0,9] = 73.595192
$row8[0,10] = 0.2952646900921413
$row8[0,11] = 0.2952646900921412
$row8[0,12] = 1.816239388880889
$row8[0,13] = 16.346154499928
$row8[0,14] = 0.01630515436136092
$row8[0,15] = 0.01630515436136092
$ws.Range("E8:T8").Value = $row8

# Row 9: FAPs -> Resolving-Mac
$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 0.07403633333333333
$row9[0,3] = 0.222109
$row9[0,4] = 0.05522216136400421
$row9[0,5] = 0.05522216136400422
$row9[0,6] = 3
$row9[0,7] = 1
$row9[0,8] = 7.746355333333334
$row9[0,9] = 23.239066
$row9[0,10] = 0.09323537902477132
$row9[0,11] = 0.0932353790247713
$row9[0,12] = 0.5735117455771112
$row9[0,13] = 5.161605710194
$row9[0,14] = 0.005148659145340015
$row9[0,15] = 0.005148659145340015
$ws.Range("E9:T9").Value = $row9

# Row 10: MuSCs -> ECs
$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 3
$row10[0,1] = 1
$row10[0,2] = 1.117936666666667
$row10[0,3] = 3.35381
$row10[0,4] = 0.8338457109086573
$row10[0,5] = 0.8338457109086574
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 3.556762333333333
$row10[0,9] = 10.670287
$row10[0,10] = 0.04280930450251701
$row10[0,11] = 0.04280930450251701
$row10[0,12] = 3.976235027052223
$row10[0,13] = 35.78611524347
$row10[0,14] = 0.03569635494640649
$row10[0,15] = 0.03569635494640648
$ws.Range("E10:T10").Value = $row10

# Row 11: MuSCs -> FAPs
$row11 = New-Object 'object[,]' 1,16
$row11[0,0] = 3
$row11[0,1] = 1
$row11[0,2] = 1.117936666666667
$row11[0,3] = 3.35381
$row11[0,4] = 0.8338457109086573
$row11[0,5] = 0.8338457109086574
$row11[0,6] = 3
$row11[0,7] = 1
$row11[0,8] = 47.24901333333333
$row11[0,9] = 141.74704
$row11[0,10] = 0.5686906263805706
$row11[0,11] = 0.5686906263805704
$row11[0,12] = 52.82140446915556
$row11[0,13] = 475.3926402224
$row11[0,14] = 0.4742002396413965
$row11[0,15] = 0.4742002396413965
$ws.Range("E11:T11").Value = $row11

# Row 12: MuSCs -> MuSCs
$row12 = New-Object 'object[,]' 1,16
$row12[0,0] = 3
$row12[0,1] = 1
$row12[0,2] = 1.117936666666667
$row12[0,3] = 3.35381
$row12[0,4] = 0.8338457109086573
$row12[0,5] = 0.8338457109086574
$row12[0,6] = 3
$row12[0,7] = 1
$row12[0,8] = 24.53173066666666
$row12[0,9] = 73.595192
$row12[0,10] = 0.2952646900921413
$row12[0,11] = 0.2952646900921412
$row12[0,12] = 27.42492120905778
$row12[0,13] = 246.82429088152
$row12[0,14] = 0.2462051954161059
$row12[0,15] = 0.2462051954161059
$ws.Range("E12:T12").Value = $row12

# Row 13: MuSCs -> Resolving-Mac
$row13 = New-Object 'object[,]' 1,16
$row13[0,0] = 3
$row13[0,1] = 1
$row13[0,2] = 1.117936666666667
$row13[0,3] = 3.35381
$row13[0,4] = 0.8338457109086573
$row13[0,5] = 0.8338457109086574
$row13[0,6] = 3
$row13[0,7] = 1
$row13[0,8] = 7.746355333333334
$row13[0,9] = 23.239066
$row13[0,10] = 0.09323537902477132
$row13[0,11] = 0.0932353790247713
$row13[0,12] = 8.659934660162223
$row13[0,13] = 77.93941194146001
$row13[0,14] = 0.07774392090474856
$row13[0,15] = 0.07774392090474856
$ws.Range("E13:T13").Value = $row13
